$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the week-ending date header
$ws.Range("F1").Value = "Saturday, January 31, 2026"

# Update the day-of-week date row (row 5, columns F:K)
$ws.Range("F5").Value = "01/26"
$ws.Range("G5").Value = "01/27"
$ws.Range("H5").Value = "01/28"
$ws.Range("I5").Value = "01/29"
$ws.Range("J5").Value = "01/30"
$ws.Range("K5").Value = "01/31"

# New job section header
$ws.Range("A7").Value = "JOB 225010 Beador Rt 15 08-1K4004"
$ws.Range("A7").Font.Bold = $true
$ws.Range("A7").Font.Underline = $true

# Crew roster rows
$ws.Range("A8").Value = "Marin"
$ws.Range("B8").Value = "Jesus"
$ws.Range("C8").Value = "LBSC1"
$ws.Range("F8").Value = "225010"
$ws.Range("G8").Value = "225010"
$ws.Range("H8").Value = "225010"

$ws.Range("A9").Value = "Palafox"
$ws.Range("B9").Value = "Sergio"
$ws.Range("C9").Value = "LBSD2"
$ws.Range("F9").Value = "225010"
$ws.Range("G9").Value = "225010"
$ws.Range("H9").Value = "225010"

$ws.Range("A10").Value = "Espinoza"
$ws.Range("B10").Value = "Luis"
$ws.Range("C10").Value = "LBSC0"
$ws.Range("G10").Value = "225010"
$ws.Range("H10").Value = "225010"

$ws.Range("A11").Value = "Richards"
$ws.Range("B11").Value = "Doug"
$ws.Range("C11").Value = "OPER1"
$ws.Range("G11").Value = "225010"
$ws.Range("H11").Value = "225010"

$ws.Range("A12").Value = "Aguirre"
$ws.Range("B12").Value = "Santiago"
$ws.Range("C12").Value = "LBSC1"
$ws.Range("F12").Value = "225010"
$ws.Range("G12").Value = "225010"
$ws.Range("H12").Value = "225010"

$ws.Range("A13").Value = "Valdivia"
$ws.Range("B13").Value = "Ivan"
$ws.Range("C13").Value = "LISC4"
$ws.Range("F13").Value = "225010"
$ws.Range("G13").Value = "225010"
$ws.Range("H13").Value = "225010"

$ws.Range("A14").Value = "Ortiz"
$ws.Range("B14").Value = "Gilberto"
$ws.Range("C14").Value = "LBSC0"
$ws.Range("F14").Value = "225010"
$ws.Range("G14").Value = "225010"
$ws.Range("H14").Value = "225010"

$ws.Range("A15").Value = "Garcia"
$ws.Range("B15").Value = "Jesus"
$ws.Range("C15").Value = "LISC6"
$ws.Range("F15").Value = "225010"
$ws.Range("G15").Value = "225010"
$ws.Range("H15").Value = "225010"

$ws.Range("A16").Value = "Avila"
$ws.Range("B16").Value = "Agustin"
$ws.Range("C16").Value = "LBSC7"
$ws.Range("F16").Value = "225010"
$ws.Range("G16").Value = "225010"
$ws.Range("H16").Value = "225010"
